$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "87÷3=29, 0" "29÷4=7, 1"
Replace-Text "18÷7=2, 4" "25÷5=5, 0"
Replace-Text "20÷2=10, 0" "50÷7=7, 1"
Replace-Text "59÷5=11, 4" "33÷6=5, 3"
Replace-Text "63÷7=9, 0" "17÷8=2, 1"
Replace-Text "58÷3=19, 1" "81÷3=27, 0"
Replace-Text "63÷9=7, 0" "45÷8=5, 5"
Replace-Text "57÷6=9, 3" "48÷7=6, 6"
Replace-Text "43÷5=8, 3" "48÷5=9, 3"
Replace-Text "45÷2=22, 1" "43÷8=5, 3"
Replace-Text "60÷6=10, 0" "79÷5=15, 4"
Replace-Text "11÷9=1, 2" "30÷6=5, 0"
Replace-Text "90÷7=12, 6" "49÷6=8, 1"
Replace-Text "62÷5=12, 2" "63÷8=7, 7"
Replace-Text "56÷3=18, 2" "86÷6=14, 2"
Replace-Text "30÷5=6, 0" "93÷6=15, 3"
Replace-Text "88÷3=29, 1" "60÷8=7, 4"
Replace-Text "88÷4=22, 0" "48÷6=8, 0"
Replace-Text "13÷2=6, 1" "16÷8=2, 0"
Replace-Text "26÷8=3, 2" "41÷4=10, 1"
Replace-Text "52÷4=13, 0" "49÷8=6, 1"
Replace-Text "57÷7=8, 1" "51÷5=10, 1"
Replace-Text "74÷7=10, 4" "59÷4=14, 3"
Replace-Text "39÷7=5, 4" "27÷6=4, 3"
Replace-Text "86÷9=9, 5" "76÷6=12, 4"
